# Update "想去人数" (column F) counts across all 4 sheets of the
# Hangzhou ACG-event workbook, per the latest scrape (commit 456a3b4).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(3, 6).Value = 110  # F3: 109 -> 110
$ws.Cells.Item(6, 6).Value = 34  # F6: 33 -> 34
$ws.Cells.Item(10, 6).Value = 10099  # F10: 10077 -> 10099
$ws.Cells.Item(11, 6).Value = 179  # F11: 178 -> 179
$ws.Cells.Item(15, 6).Value = 1948  # F15: 1946 -> 1948
$ws.Cells.Item(20, 6).Value = 153  # F20: 152 -> 153
$ws.Cells.Item(22, 6).Value = 216  # F22: 215 -> 216
$ws.Cells.Item(23, 6).Value = 1099  # F23: 1097 -> 1099
$ws.Cells.Item(24, 6).Value = 73  # F24: 70 -> 73
$ws.Cells.Item(26, 6).Value = 613  # F26: 610 -> 613
$ws.Cells.Item(28, 6).Value = 138  # F28: 136 -> 138
$ws.Cells.Item(29, 6).Value = 614  # F29: 610 -> 614
$ws.Cells.Item(30, 6).Value = 2709  # F30: 2694 -> 2709
$ws.Cells.Item(32, 6).Value = 621  # F32: 612 -> 621
$ws.Cells.Item(35, 6).Value = 15  # F35: 14 -> 15
$ws.Cells.Item(36, 6).Value = 494  # F36: 486 -> 494
$ws.Cells.Item(37, 6).Value = 202  # F37: 201 -> 202
$ws.Cells.Item(38, 6).Value = 13  # F38: 12 -> 13
$ws.Cells.Item(39, 6).Value = 1199  # F39: 1193 -> 1199
$ws.Cells.Item(40, 6).Value = 16  # F40: 2 -> 16
$ws.Cells.Item(41, 6).Value = 207  # F41: 199 -> 207
$ws.Cells.Item(42, 6).Value = 103  # F42: 102 -> 103
$ws.Cells.Item(43, 6).Value = 58  # F43: 57 -> 58
$ws.Cells.Item(44, 6).Value = 104  # F44: 102 -> 104
$ws.Cells.Item(45, 6).Value = 133  # F45: 120 -> 133
$ws.Cells.Item(46, 6).Value = 29  # F46: 28 -> 29
$ws.Cells.Item(47, 6).Value = 4039  # F47: 4035 -> 4039

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(5, 6).Value = 4036  # F5: 4035 -> 4036
$ws.Cells.Item(7, 6).Value = 42  # F7: 39 -> 42

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 728  # F2: 727 -> 728
$ws.Cells.Item(3, 6).Value = 384  # F3: 382 -> 384
$ws.Cells.Item(4, 6).Value = 37  # F4: 36 -> 37

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 728  # F2: 727 -> 728
$ws.Cells.Item(3, 6).Value = 384  # F3: 382 -> 384
$ws.Cells.Item(7, 6).Value = 110  # F7: 109 -> 110
$ws.Cells.Item(8, 6).Value = 37  # F8: 36 -> 37
$ws.Cells.Item(9, 6).Value = 34  # F9: 33 -> 34
$ws.Cells.Item(11, 6).Value = 69  # F11: 68 -> 69
$ws.Cells.Item(12, 6).Value = 10099  # F12: 10077 -> 10099
$ws.Cells.Item(13, 6).Value = 179  # F13: 178 -> 179
$ws.Cells.Item(16, 6).Value = 1948  # F16: 1946 -> 1948
$ws.Cells.Item(21, 6).Value = 216  # F21: 215 -> 216
$ws.Cells.Item(22, 6).Value = 1099  # F22: 1097 -> 1099
$ws.Cells.Item(23, 6).Value = 73  # F23: 70 -> 73
$ws.Cells.Item(25, 6).Value = 4036  # F25: 4035 -> 4036
$ws.Cells.Item(27, 6).Value = 613  # F27: 610 -> 613
$ws.Cells.Item(29, 6).Value = 138  # F29: 136 -> 138
$ws.Cells.Item(30, 6).Value = 614  # F30: 610 -> 614
$ws.Cells.Item(31, 6).Value = 2710  # F31: 2694 -> 2710
$ws.Cells.Item(33, 6).Value = 42  # F33: 39 -> 42
$ws.Cells.Item(35, 6).Value = 621  # F35: 612 -> 621
$ws.Cells.Item(37, 6).Value = 494  # F37: 486 -> 494
$ws.Cells.Item(39, 6).Value = 202  # F39: 201 -> 202
$ws.Cells.Item(40, 6).Value = 58  # F40: 57 -> 58
$ws.Cells.Item(41, 6).Value = 104  # F41: 102 -> 104
$ws.Cells.Item(42, 6).Value = 133  # F42: 120 -> 133
$ws.Cells.Item(43, 6).Value = 29  # F43: 28 -> 29
$ws.Cells.Item(44, 6).Value = 4039  # F44: 4035 -> 4039
